$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 261
$ws.Range("C3").Value = 159619
$ws.Range("C4").Value = 150665
$ws.Range("C8").Value = 64.14
